# Updated experiment result with new k for cross-validation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1 was blank, now labelled "Row"
$ws.Range("A1").Value = "Row"

# Row 2 (random_forest) - refreshed metrics
$ws.Range("B2").Value = 3.4015160345908879
$ws.Range("C2").Value = 0.24073008029659504
$ws.Range("D2").Value = 2.6230186965811959
$ws.Range("E2").Value = 0.36416577228382507
$ws.Range("F2").Value = 0.60346149196433829
$ws.Range("G2").Value = 0.63266249314548861
$ws.Range("H2").Value = 0.63583422771617493
$ws.Range("I2").Value = 0.79840499447713709

# Row 3 (lsboost) - refreshed metrics
$ws.Range("B3").Value = 3.6665942560691902
$ws.Range("C3").Value = 0.25949003935379972
$ws.Range("D3").Value = 3.1032471973291464
$ws.Range("E3").Value = 0.42313581534237832
$ws.Range("F3").Value = 0.65048890485724531
$ws.Range("G3").Value = 0.74849184691971693
$ws.Range("H3").Value = 0.57686418465762168
$ws.Range("I3").Value = 0.78356355667974831

# Row 4 (old_model) values are unchanged by this commit.

# Column widths tightened up after the refresh
$ws.Columns.Item(1).ColumnWidth = 17.307291666666668
$ws.Columns.Item(2).ColumnWidth = 7.022135416666667
$ws.Columns.Item(3).ColumnWidth = 8.877604166666666
$ws.Columns.Item(4).ColumnWidth = 6.877604166666667
$ws.Columns.Item(5).ColumnWidth = 6.877604166666667
$ws.Columns.Item(6).ColumnWidth = 6.877604166666667
$ws.Columns.Item(7).ColumnWidth = 6.877604166666667
$ws.Columns.Item(8).ColumnWidth = 6.877604166666667
$ws.Columns.Item(9).ColumnWidth = 11.877604166666666
